$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "568.65", "1.00") are preserved verbatim instead of being
# converted into floating point numbers that lose formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.940.89"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "2.613.29"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "568.65"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "143.02"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "2.639.46"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "0.374"
$ws.Range("E13").Value = "  +8.11%  "
$ws.Range("D14").Value = "3.074.72"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "60.859.03"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "23.54"
$ws.Range("E16").Value = "  +4.36%  "
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "2.621.77"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +3.31%  "
$ws.Range("D20").Value = "349.56"
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("D21").Value = "10.97"
$ws.Range("E21").Value = "  +7.11%  "
$ws.Range("D22").Value = "7.01"
$ws.Range("E22").Value = "  +12.84%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +14.21%  "
$ws.Range("D25").Value = "63.55"
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "7.73"
$ws.Range("E28").Value = "  +6.18%  "
$ws.Range("D29").Value = "0.0₃0795"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").Value = "1.82"
$ws.Range("E30").Value = "  +7.67%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "6.28"
$ws.Range("E32").Value = "  +3.11%  "
$ws.Range("D33").Value = "161.59"
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("D34").Value = "19.57"
$ws.Range("E34").Value = "  +2.54%  "
$ws.Range("D35").Value = "4.22"
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("D36").Value = "0.969"
$ws.Range("E36").Value = "  +9.69%  "
$ws.Range("E37").Value = "  +4.74%  "
$ws.Range("D38").Value = "1.60"
$ws.Range("E38").Value = "  +6.46%  "
$ws.Range("D39").Value = "37.71"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("D40").Value = "0.856"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").Value = "3.81"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("D42").Value = "303.70"
$ws.Range("E42").Value = "  +3.10%  "
$ws.Range("D43").Value = "140.52"
$ws.Range("E43").Value = "  +12.96%  "
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "0.0986"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("D46").Value = "0.607"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").Value = "0.0549"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("E48").Value = "  +4.01%  "
$ws.Range("D49").Value = "10.69"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").Value = "4.83"
$ws.Range("E50").Value = "  +7.12%  "
$ws.Range("D51").Value = "19.49"
$ws.Range("E51").Value = "  +5.20%  "

# Restore the default (General/Normal) style on column D now that the
# text values are committed, matching the original workbook styling.
$ws.Range("D2:D51").Style = "Normal"

